$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.645526885986328
$ws.Range("B1").Value = 2.436991691589355
$ws.Range("C1").Value = 2.691607475280762
$ws.Range("D1").Value = 3.005879163742065
$ws.Range("E1").Value = 0.4752430617809296
